$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Review")

$ws.Range("M6").Value = "no"
$ws.Range("M7").Value = "no"
$ws.Range("M8").Value = "no"
$ws.Range("M9").Value = "no"

$ws.Select()
$ws.Range("O9").Select()
